$d = $word.ActiveDocument

# 1. Locate the "Authors" paragraph, currently reading
#    " Alexandros Konstantakos, Zahid Ibnu Yusuf".
$authorsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Alexandros Konstantakos, Zahid Ibnu Yusuf") {
        $authorsPara = $p
        break
    }
}

# 2. Trim it down to just Alexandros' name plus his matriculation number
#    (assigning straight to Range.Text keeps the run's existing formatting
#    / identity instead of fabricating a brand-new run).
$authorsPara.Range.Text = " Alexandros Konstantakos, 741590"

# 3. Insert a brand-new paragraph right after it for the second author; it
#    automatically inherits the centered / Ebrima / bold / 24pt formatting
#    used throughout this authors block.
$authorsPara.Range.InsertParagraphAfter()
$newPara = $authorsPara.Next()

# 4. Fill the new (still empty) paragraph with two separate runs --
#    "Zahid Ibnu Yusuf" and ", 741463" -- via a raw OOXML fragment so the
#    run split is explicit rather than left to auto-merging.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Ebrima" w:hAnsi="Ebrima"/><w:b/><w:sz w:val="24"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Ebrima" w:hAnsi="Ebrima"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>Zahid Ibnu Yusuf</w:t></w:r>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Ebrima" w:hAnsi="Ebrima"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>, 741463</w:t></w:r>' +
       '</w:p>' +
       '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $newPara.Range.InsertXML($xml)
